$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '45.344.53'
$cell.Style = $origStyle
$ws.Range("E2").Value = '  -0.23%  '

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.369.75'
$cell.Style = $origStyle
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("E4").Value = '  +0.07%  '

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '313.69'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  -1.24%  '

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '107.94'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  -3.64%  '

$ws.Range("E7").Value = '  -0.89%  '

$ws.Range("E8").Value = '  +0.07%  '

$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.613'
$cell.Style = $origStyle
$ws.Range("E9").Value = '  -2.92%  '

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '40.78'
$cell.Style = $origStyle
$ws.Range("E10").Value = '  -3.59%  '

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0917'
$cell.Style = $origStyle
$ws.Range("E11").Value = '  -1.52%  '

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.48'
$cell.Style = $origStyle
$ws.Range("E12").Value = '  -2.30%  '

$ws.Range("E13").Value = '  +1.09%  '

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.981'
$cell.Style = $origStyle
$ws.Range("E14").Value = '  -3.70%  '

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.734.17'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  -0.36%  '

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.34'
$cell.Style = $origStyle
$ws.Range("E16").Value = '  -3.01%  '

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.349.41'
$cell.Style = $origStyle
$ws.Range("E17").Value = '  -1.53%  '

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '45.357.00'
$cell.Style = $origStyle
$ws.Range("E18").Value = '  -0.18%  '

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.55'
$cell.Style = $origStyle
$ws.Range("E19").Value = '  +18.79%  '

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.28'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  -5.19%  '

$ws.Range("E21").Value = '  -2.18%  '

$ws.Range("B22").Value = 'PancakeSwap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.60'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '73.36'
$cell.Style = $origStyle
$ws.Range("E23").Value = '  -2.32%  '

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '261.04'
$cell.Style = $origStyle
$ws.Range("E24").Value = '  -3.68%  '

$ws.Range("E25").Value = '  +0.87%  '

$ws.Range("E26").Value = '  +0.13%  '

$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.52'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  -0.15%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '11.13'
$cell.Style = $origStyle
$ws.Range("E28").Value = '  -1.41%  '

$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.Style = $origStyle
$ws.Range("E29").Value = '  -1.72%  '

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0966'
$cell.Style = $origStyle
$ws.Range("E30").Value = '  +2.10%  '

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '22.30'
$cell.Style = $origStyle
$ws.Range("E31").Value = '  -3.02%  '

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '37.00'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  -4.27%  '

$ws.Range("E33").Value = '  -1.94%  '

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.90'
$cell.Style = $origStyle
$ws.Range("E34").Value = '  -3.38%  '

$ws.Range("E35").Value = '  -2.15%  '

$ws.Range("E36").Value = '  -0.83%  '

$ws.Range("E37").Value = '  -4.29%  '

$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.89'
$cell.Style = $origStyle
$ws.Range("E38").Value = '  +8.38%  '

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '3.94'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  +0.11%  '

$ws.Range("E40").Value = '  -4.48%  '

$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0354'
$cell.Style = $origStyle
$ws.Range("E41").Value = '  -2.93%  '

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '98.64'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  -6.41%  '

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '70.08'
$cell.Style = $origStyle
$ws.Range("E43").Value = '  -1.70%  '

$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '13.09'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("E45").Value = '  -5.96%  '

$ws.Range("E46").Value = '  -0.02%  '

$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.94'
$cell.Style = $origStyle
$ws.Range("E47").Value = '  +2.27%  '

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.814.79'
$cell.Style = $origStyle
$ws.Range("E48").Value = '  +9.83%  '

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '82.60'
$cell.Style = $origStyle
$ws.Range("E49").Value = '  +4.40%  '

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '111.14'
$cell.Style = $origStyle
$ws.Range("E50").Value = '  -6.15%  '

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '9.19'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  -1.32%  '
